$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# "SHOW SCHEDULES [LIKE 'likeName']" is removed from the command set; the
# row that used to show it now documents "LIST SCHEDULES [LIKE 'likeName']"
# instead (duplicate command rationalized away).
$ws.Range("A14").Value = "LIST SCHEDULES [LIKE 'likeName']"

# The row that used to document "LIST SCHEDULES [LIKE 'likeName']" is
# repurposed as a not-yet-implemented placeholder. The leading apostrophe
# forces Excel to store it as literal text (quote-prefixed) since it starts
# with "-", matching how Excel marks such cells internally.
$ws.Range("A16").Value = "'-- not implemented --"
